$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.975.71"
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("D3").Value = "3.031.26"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "593.64"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "153.81"
$ws.Range("E6").Value = "  +7.97%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.028.07"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +14.51%  "
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").Value = "35.68"
$ws.Range("E14").Value = "  +4.95%  "
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "3.532.91"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("D18").Value = "62.903.55"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("D19").Value = "3.032.51"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "452.32"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "0.697"
$ws.Range("E22").Value = "  +2.97%  "
$ws.Range("D23").Value = "7.50"
$ws.Range("E23").Value = "  +3.49%  "
$ws.Range("D24").Value = "83.08"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "11.29"
$ws.Range("E25").Value = "  +9.20%  "
$ws.Range("E26").Value = "  +6.13%  "
$ws.Range("E27").Value = "  +5.19%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "7.47"
$ws.Range("E29").Value = "  +5.63%  "
$ws.Range("E30").Value = "  +10.89%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "27.55"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").Value = "0.0₃0859"
$ws.Range("E35").Value = "  +6.37%  "
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("D37").Value = "5.91"
$ws.Range("E37").Value = "  +3.52%  "
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  +12.24%  "
$ws.Range("E39").Value = "  +7.93%  "
$ws.Range("D40").Value = "2.09"
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").Value = "50.52"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "0.307"
$ws.Range("E43").Value = "  +15.98%  "
$ws.Range("D44").Value = "43.72"
$ws.Range("E44").Value = "  +11.95%  "
$ws.Range("D45").Value = "390.02"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "0.0360"
$ws.Range("E46").Value = "  +3.75%  "
$ws.Range("D47").Value = "2.721.40"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").Value = "133.54"
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("E50").Value = "  +7.87%  "
$ws.Range("D51").Value = "24.98"
$ws.Range("E51").Value = "  +8.22%  "